$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1 (paragraph "Nel primo caso, il programma chiedera...")
#   Split the run so that the space between "massimo" and "di 8" becomes
#   its own bold run, and a "_GoBack" bookmark is inserted between
#   "Succe" and "ssivamente".
# ---------------------------------------------------------------------

# 1a. Make the single space after "massimo" bold (this naturally splits
#     the run into "...un massimo" | " " (bold) | "di 8 giocatori...").
$rng1 = $d.Content
$rng1.Find.Execute("massimo")
$spaceRng = $d.Range($rng1.End, $rng1.End + 1)
$spaceRng.Bold = 1

# 1b. Insert the "_GoBack" bookmark right after "Succe" (splitting that
#     run into "...di 8 giocatori. Succe" | "ssivamente, ...").
$rng2 = $d.Content
$rng2.Find.Execute("Succe")
$bmRng = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------
# Edit 2 (paragraph "L'esito della partita verra deciso...")
#   Remove one of the two spaces between "partita." and "Verranno", and
#   split the run into three pieces at "che decr|eteranno" and at
#   "esito della partita.| Verranno".
# ---------------------------------------------------------------------

# The original run ends right before "vincente, 0 punto al perdente...",
# which is (and must stay) its own separate, untouched run.
$origRunEnd = $d.Content
$origRunEnd.Find.Execute("giocatore vincente")
$runBoundary = $origRunEnd.Start + "giocatore ".Length

# 2a. Collapse the double space before "Verranno" down to a single space.
$rng3 = $d.Content
$rng3.Find.Execute("esito della partita.")
$extraSpace = $d.Range($rng3.End, $rng3.End + 1)
$extraSpace.Text = ""
$runBoundary = $runBoundary - 1

# 2b. Split right after "che decr" (-> "...che decr" | "eteranno...").
$rng4 = $d.Content
$rng4.Find.Execute("che decr")
$afterDecr = $d.Range($rng4.End, $runBoundary)
$afterDecr.Bold = 1
$afterDecr.Bold = 0

# 2c. Split right after "esito della partita." (-> "...partita." | " Verranno...").
$rng5 = $d.Content
$rng5.Find.Execute("esito della partita.")
$afterPartita = $d.Range($rng5.End, $runBoundary)
$afterPartita.Bold = 1
$afterPartita.Bold = 0

Write-Output "Paragraph 4: $($d.Paragraphs(4).Range.Text)"
Write-Output "Paragraph 13: $($d.Paragraphs(13).Range.Text)"
